$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 1827.5714
$ws.Range("I28").Value = 1827.5714
$ws.Range("K28").Value = 1827.5714
$ws.Range("M28").Value = -1342.5714

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1999
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 3450
$ws.Range("J51").Value = 3450
$ws.Range("L51").Value = 3450
$ws.Range("N51").Value = -4418

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1296.8
$ws.Range("I80").Value = 1498
$ws.Range("J80").Value = 995
$ws.Range("K80").Value = 4494
$ws.Range("L80").Value = 2985
$ws.Range("M80").Value = -3496
$ws.Range("N80").Value = -4981

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 1296.8
$ws.Range("I83").Value = 1498
$ws.Range("J83").Value = 995
$ws.Range("K83").Value = 13482
$ws.Range("L83").Value = 8955
$ws.Range("M83").Value = -8490
$ws.Range("N83").Value = -18939

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 9195.799999999999
$ws.Range("J125").Value = 8662.333000000001
$ws.Range("L125").Value = 77960.997
$ws.Range("N125").Value = -82880.997

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 2000
$ws.Range("I129").Value = 2000
$ws.Range("K129").Value = 6000
$ws.Range("M129").Value = -1000

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2882.9333
$ws.Range("I132").Value = 1921.5
$ws.Range("K132").Value = 5764.5
$ws.Range("M132").Value = -3234.5

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1412
$ws.Range("I137").Value = 1412
$ws.Range("K137").Value = 4236
$ws.Range("M137").Value = -1686

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 6915.8335
$ws.Range("I141").Value = 6299
$ws.Range("K141").Value = 18897
$ws.Range("M141").Value = -13717

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3086.125
$ws.Range("I45").Value = 3086.125
$ws.Range("K45").Value = 3086.125
$ws.Range("M45").Value = -2709.125

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 457.66666
$ws.Range("I22").Value = 457.66666
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 457.66666
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -284.66666
$ws.Range("N22").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1571.5
$ws.Range("I105").Value = 1571.5
$ws.Range("K105").Value = 1571.5
$ws.Range("M105").Value = 175.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2246.1
$ws.Range("I107").Value = 2508.8572
$ws.Range("K107").Value = 2508.8572
$ws.Range("M107").Value = -588.8571999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 125000
$ws.Range("J132").Value = 125000
$ws.Range("L132").Value = 125000
$ws.Range("N132").Value = -135120

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4017.2
$ws.Range("I99").Value = 4017.2
$ws.Range("K99").Value = 4017.2
$ws.Range("M99").Value = -2519.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4017.2
$ws.Range("I126").Value = 4017.2
$ws.Range("K126").Value = 12051.6
$ws.Range("M126").Value = -9581.599999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 4231.3335
$ws.Range("J132").Value = 5122.75
$ws.Range("L132").Value = 15368.25
$ws.Range("N132").Value = -20428.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 62.5
$ws.Range("J2").Value = 100
$ws.Range("L2").Value = 600
$ws.Range("N2").Value = -826

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 929.5
$ws.Range("J12").Value = 922
$ws.Range("L12").Value = 2766
$ws.Range("N12").Value = -3112

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 1720
$ws.Range("J22").Value = 1720
$ws.Range("L22").Value = 5160
$ws.Range("N22").Value = -5498

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H27").Value = 1720
$ws.Range("J27").Value = 1720
$ws.Range("L27").Value = 5160
$ws.Range("N27").Value = -5364

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H38").Value = 309.16666
$ws.Range("J38").Value = 226
$ws.Range("L38").Value = 678
$ws.Range("N38").Value = -1372

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 2059
$ws.Range("J107").Value = 721.3333
$ws.Range("L107").Value = 2163.9999
$ws.Range("N107").Value = -6003.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 900
$ws.Range("I108").Value = 900
$ws.Range("K108").Value = 2700
$ws.Range("M108").Value = 180

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 7361.1113

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 4952.2856
$ws.Range("I116").Value = 3998
$ws.Range("J116").Value = 5000
$ws.Range("K116").Value = 11994
$ws.Range("L116").Value = 15000
$ws.Range("M116").Value = -8552
$ws.Range("N116").Value = -21884

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 1392.625
$ws.Range("I117").Value = 2196.6667
$ws.Range("J117").Value = 910.2
$ws.Range("K117").Value = 6590.000100000001
$ws.Range("L117").Value = 2730.6
$ws.Range("M117").Value = -3148.000100000001
$ws.Range("N117").Value = -9614.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H118").Value = 5000
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 5000
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 15000
$ws.Range("M118").ClearContents()
$ws.Range("N118").Value = -17486

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 1286.1818
$ws.Range("J121").Value = 1541.4286
$ws.Range("L121").Value = 4624.2858
$ws.Range("N121").Value = -7244.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 1121.25
$ws.Range("I129").Value = 993.5
$ws.Range("J129").Value = 1249
$ws.Range("K129").Value = 2980.5
$ws.Range("L129").Value = 3747
$ws.Range("M129").Value = 2019.5
$ws.Range("N129").Value = -13747

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6524.1665
$ws.Range("I80").Value = 6126.4287
$ws.Range("K80").Value = 6126.4287
$ws.Range("M80").Value = -5128.4287

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 6524.1665
$ws.Range("I83").Value = 6126.4287
$ws.Range("K83").Value = 30632.1435
$ws.Range("M83").Value = -25640.1435

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1216.8667
$ws.Range("I16").Value = 1061
$ws.Range("K16").Value = 1061
$ws.Range("M16").Value = -891

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 190.4
$ws.Range("I55").Value = 188.66667
$ws.Range("K55").Value = 188.66667
$ws.Range("M55").Value = -15.66667000000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2598.5293
$ws.Range("I132").Value = 2450
$ws.Range("J132").Value = 2955
$ws.Range("K132").Value = 7350
$ws.Range("L132").Value = 8865
$ws.Range("M132").Value = -4820
$ws.Range("N132").Value = -13925

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 2972.45
$ws.Range("I132").Value = 2036
$ws.Range("J132").Value = 3596.75
$ws.Range("K132").Value = 6108
$ws.Range("L132").Value = 10790.25
$ws.Range("M132").Value = -3578
$ws.Range("N132").Value = -15850.25
